$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 7).Value2 = 0.340305
$ws.Cells.Item(2, 8).Value2 = 1.020915
$ws.Cells.Item(2, 9).Value2 = 0.2413792532744959
$ws.Cells.Item(2, 10).Value2 = 0.2522899505114672
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 13).Value2 = 1.919165666666667
$ws.Cells.Item(2, 14).Value2 = 5.757497
$ws.Cells.Item(2, 15).Value2 = 0.09179925534063474
$ws.Cells.Item(2, 16).Value2 = 0.1039451703609422
$ws.Cells.Item(2, 17).Value2 = 0.653101672195
$ws.Cells.Item(2, 18).Value2 = 5.877915049755
$ws.Cells.Item(2, 19).Value2 = 0.02215843570527719
$ws.Cells.Item(2, 20).Value2 = 0.02622432188626814

$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 7).Value2 = 0.340305
$ws.Cells.Item(3, 8).Value2 = 1.020915
$ws.Cells.Item(3, 9).Value2 = 0.2413792532744959
$ws.Cells.Item(3, 10).Value2 = 0.2522899505114672
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 13).Value2 = 8.648731
$ws.Cells.Item(3, 14).Value2 = 25.946193
$ws.Cells.Item(3, 15).Value2 = 0.4136938666792861
$ws.Cells.Item(3, 16).Value2 = 0.4684295018482661
$ws.Cells.Item(3, 17).Value2 = 2.943206402955
$ws.Cells.Item(3, 18).Value2 = 26.488857626595
$ws.Cells.Item(3, 19).Value2 = 0.09985711662328495
$ws.Cells.Item(3, 20).Value2 = 0.1181800558394103

$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 7).Value2 = 0.340305
$ws.Cells.Item(4, 8).Value2 = 1.020915
$ws.Cells.Item(4, 9).Value2 = 0.2413792532744959
$ws.Cells.Item(4, 10).Value2 = 0.2522899505114672
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 13).Value2 = 1.069144
$ws.Cells.Item(4, 14).Value2 = 3.207432
$ws.Cells.Item(4, 15).Value2 = 0.05114025576665047
$ws.Cells.Item(4, 16).Value2 = 0.05790659824245461
$ws.Cells.Item(4, 17).Value2 = 0.3638350489200001
$ws.Cells.Item(4, 18).Value2 = 3.274515440280001
$ws.Cells.Item(4, 19).Value2 = 0.01234419674922083
$ws.Cells.Item(4, 20).Value2 = 0.01460925280487629

$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 7).Value2 = 0.340305
$ws.Cells.Item(5, 8).Value2 = 1.020915
$ws.Cells.Item(5, 9).Value2 = 0.2413792532744959
$ws.Cells.Item(5, 10).Value2 = 0.2522899505114672
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 13).Value2 = 1.940482
$ws.Cells.Item(5, 14).Value2 = 5.821446
$ws.Cells.Item(5, 15).Value2 = 0.09281887733605711
$ws.Cells.Item(5, 16).Value2 = 0.1050996980488267
$ws.Cells.Item(5, 17).Value2 = 0.6603557270100001
$ws.Cells.Item(5, 18).Value2 = 5.94320154309
$ws.Cells.Item(5, 19).Value2 = 0.0224045513011545
$ws.Cells.Item(5, 20).Value2 = 0.02651559761950864

$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 7).Value2 = 0.340305
$ws.Cells.Item(6, 8).Value2 = 1.020915
$ws.Cells.Item(6, 9).Value2 = 0.2413792532744959
$ws.Cells.Item(6, 10).Value2 = 0.2522899505114672
$ws.Cells.Item(6, 11).Value2 = 2
$ws.Cells.Item(6, 13).Value2 = 7.328590999999999
$ws.Cells.Item(6, 14).Value2 = 14.657182
$ws.Cells.Item(6, 15).Value2 = 0.3505477448773717
$ws.Cells.Item(6, 16).Value2 = 0.2646190314995103
$ws.Cells.Item(6, 17).Value2 = 2.493956160255
$ws.Cells.Item(6, 18).Value2 = 14.96373696153
$ws.Cells.Item(6, 19).Value2 = 0.08461495289555848
$ws.Cells.Item(6, 20).Value2 = 0.06676072236140382

$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 7).Value2 = 0.472982
$ws.Cells.Item(7, 8).Value2 = 1.418946
$ws.Cells.Item(7, 9).Value2 = 0.3354874068035369
$ws.Cells.Item(7, 10).Value2 = 0.3506519309819567
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 13).Value2 = 1.919165666666667
$ws.Cells.Item(7, 14).Value2 = 5.757497
$ws.Cells.Item(7, 15).Value2 = 0.09179925534063474
$ws.Cells.Item(7, 16).Value2 = 0.1039451703609422
$ws.Cells.Item(7, 17).Value2 = 0.9077308153513333
$ws.Cells.Item(7, 18).Value2 = 8.169577338162
$ws.Cells.Item(7, 19).Value2 = 0.03079749412072528
$ws.Cells.Item(7, 20).Value2 = 0.03644857470331284

$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 7).Value2 = 0.472982
$ws.Cells.Item(8, 8).Value2 = 1.418946
$ws.Cells.Item(8, 9).Value2 = 0.3354874068035369
$ws.Cells.Item(8, 10).Value2 = 0.3506519309819567
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 13).Value2 = 8.648731
$ws.Cells.Item(8, 14).Value2 = 25.946193
$ws.Cells.Item(8, 15).Value2 = 0.4136938666792861
$ws.Cells.Item(8, 16).Value2 = 0.4684295018482661
$ws.Cells.Item(8, 17).Value2 = 4.090694085842
$ws.Cells.Item(8, 18).Value2 = 36.816246772578
$ws.Cells.Item(8, 19).Value2 = 0.1387890825427618
$ws.Cells.Item(8, 20).Value2 = 0.1642557093520106

$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 7).Value2 = 0.472982
$ws.Cells.Item(9, 8).Value2 = 1.418946
$ws.Cells.Item(9, 9).Value2 = 0.3354874068035369
$ws.Cells.Item(9, 10).Value2 = 0.3506519309819567
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 13).Value2 = 1.069144
$ws.Cells.Item(9, 14).Value2 = 3.207432
$ws.Cells.Item(9, 15).Value2 = 0.05114025576665047
$ws.Cells.Item(9, 16).Value2 = 0.05790659824245461
$ws.Cells.Item(9, 17).Value2 = 0.505685867408
$ws.Cells.Item(9, 18).Value2 = 4.551172806672001
$ws.Cells.Item(9, 19).Value2 = 0.01715691179042319
$ws.Cells.Item(9, 20).Value2 = 0.02030506049031309

$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 7).Value2 = 0.472982
$ws.Cells.Item(10, 8).Value2 = 1.418946
$ws.Cells.Item(10, 9).Value2 = 0.3354874068035369
$ws.Cells.Item(10, 10).Value2 = 0.3506519309819567
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 13).Value2 = 1.940482
$ws.Cells.Item(10, 14).Value2 = 5.821446
$ws.Cells.Item(10, 15).Value2 = 0.09281887733605711
$ws.Cells.Item(10, 16).Value2 = 0.1050996980488267
$ws.Cells.Item(10, 17).Value2 = 0.9178130573240001
$ws.Cells.Item(10, 18).Value2 = 8.260317515916
$ws.Cells.Item(10, 19).Value2 = 0.03113956445988938
$ws.Cells.Item(10, 20).Value2 = 0.03685341206644167

$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 7).Value2 = 0.472982
$ws.Cells.Item(11, 8).Value2 = 1.418946
$ws.Cells.Item(11, 9).Value2 = 0.3354874068035369
$ws.Cells.Item(11, 10).Value2 = 0.3506519309819567
$ws.Cells.Item(11, 11).Value2 = 2
$ws.Cells.Item(11, 13).Value2 = 7.328590999999999
$ws.Cells.Item(11, 14).Value2 = 14.657182
$ws.Cells.Item(11, 15).Value2 = 0.3505477448773717
$ws.Cells.Item(11, 16).Value2 = 0.2646190314995103
$ws.Cells.Item(11, 17).Value2 = 3.466291628362
$ws.Cells.Item(11, 18).Value2 = 20.797749770172
$ws.Cells.Item(11, 19).Value2 = 0.1176043538897373
$ws.Cells.Item(11, 20).Value2 = 0.09278917436987849

$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 7).Value2 = 0.4136363333333333
$ws.Cells.Item(12, 8).Value2 = 1.240909
$ws.Cells.Item(12, 9).Value2 = 0.2933933655608953
$ws.Cells.Item(12, 10).Value2 = 0.3066551771687498
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 13).Value2 = 1.919165666666667
$ws.Cells.Item(12, 14).Value2 = 5.757497
$ws.Cells.Item(12, 15).Value2 = 0.09179925534063474
$ws.Cells.Item(12, 16).Value2 = 0.1039451703609422
$ws.Cells.Item(12, 17).Value2 = 0.7938366494192222
$ws.Cells.Item(12, 18).Value2 = 7.144529844773
$ws.Cells.Item(12, 19).Value2 = 0.02693329248037282
$ws.Cells.Item(12, 20).Value2 = 0.03187532463287062

$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 7).Value2 = 0.4136363333333333
$ws.Cells.Item(13, 8).Value2 = 1.240909
$ws.Cells.Item(13, 9).Value2 = 0.2933933655608953
$ws.Cells.Item(13, 10).Value2 = 0.3066551771687498
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 13).Value2 = 8.648731
$ws.Cells.Item(13, 14).Value2 = 25.946193
$ws.Cells.Item(13, 15).Value2 = 0.4136938666792861
$ws.Cells.Item(13, 16).Value2 = 0.4684295018482661
$ws.Cells.Item(13, 17).Value2 = 3.577429378826333
$ws.Cells.Item(13, 18).Value2 = 32.196864409437
$ws.Cells.Item(13, 19).Value2 = 0.1213750358569361
$ws.Cells.Item(13, 20).Value2 = 0.1436463318803493

$ws.Cells.Item(14, 5).Value2 = 3
$ws.Cells.Item(14, 7).Value2 = 0.4136363333333333
$ws.Cells.Item(14, 8).Value2 = 1.240909
$ws.Cells.Item(14, 9).Value2 = 0.2933933655608953
$ws.Cells.Item(14, 10).Value2 = 0.3066551771687498
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 13).Value2 = 1.069144
$ws.Cells.Item(14, 14).Value2 = 3.207432
$ws.Cells.Item(14, 15).Value2 = 0.05114025576665047
$ws.Cells.Item(14, 16).Value2 = 0.05790659824245461
$ws.Cells.Item(14, 17).Value2 = 0.4422368039653334
$ws.Cells.Item(14, 18).Value2 = 3.980131235688
$ws.Cells.Item(14, 19).Value2 = 0.01500421175502257
$ws.Cells.Item(14, 20).Value2 = 0.01775735814327953

$ws.Cells.Item(15, 5).Value2 = 3
$ws.Cells.Item(15, 7).Value2 = 0.4136363333333333
$ws.Cells.Item(15, 8).Value2 = 1.240909
$ws.Cells.Item(15, 9).Value2 = 0.2933933655608953
$ws.Cells.Item(15, 10).Value2 = 0.3066551771687498
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 13).Value2 = 1.940482
$ws.Cells.Item(15, 14).Value2 = 5.821446
$ws.Cells.Item(15, 15).Value2 = 0.09281887733605711
$ws.Cells.Item(15, 16).Value2 = 0.1050996980488267
$ws.Cells.Item(15, 17).Value2 = 0.8026538593793333
$ws.Cells.Item(15, 18).Value2 = 7.223884734414
$ws.Cells.Item(15, 19).Value2 = 0.0272324428092097
$ws.Cells.Item(15, 20).Value2 = 0.03222936652554508

$ws.Cells.Item(16, 5).Value2 = 3
$ws.Cells.Item(16, 7).Value2 = 0.4136363333333333
$ws.Cells.Item(16, 8).Value2 = 1.240909
$ws.Cells.Item(16, 9).Value2 = 0.2933933655608953
$ws.Cells.Item(16, 10).Value2 = 0.3066551771687498
$ws.Cells.Item(16, 11).Value2 = 2
$ws.Cells.Item(16, 13).Value2 = 7.328590999999999
$ws.Cells.Item(16, 14).Value2 = 14.657182
$ws.Cells.Item(16, 15).Value2 = 0.3505477448773717
$ws.Cells.Item(16, 16).Value2 = 0.2646190314995103
$ws.Cells.Item(16, 17).Value2 = 3.031371509739667
$ws.Cells.Item(16, 18).Value2 = 18.188229058438
$ws.Cells.Item(16, 19).Value2 = 0.1028483826593542
$ws.Cells.Item(16, 20).Value2 = 0.08114679598670532

$ws.Cells.Item(17, 5).Value2 = 2
$ws.Cells.Item(17, 7).Value2 = 0.182912
$ws.Cells.Item(17, 8).Value2 = 0.365824
$ws.Cells.Item(17, 9).Value2 = 0.129739974361072
$ws.Cells.Item(17, 10).Value2 = 0.09040294133782634
$ws.Cells.Item(17, 11).Value2 = 3
$ws.Cells.Item(17, 13).Value2 = 1.919165666666667
$ws.Cells.Item(17, 14).Value2 = 5.757497
$ws.Cells.Item(17, 15).Value2 = 0.09179925534063474
$ws.Cells.Item(17, 16).Value2 = 0.1039451703609422
$ws.Cells.Item(17, 17).Value2 = 0.3510384304213334
$ws.Cells.Item(17, 18).Value2 = 2.106230582528
$ws.Cells.Item(17, 19).Value2 = 0.01191003303425945
$ws.Cells.Item(17, 20).Value2 = 0.009396949138490625

$ws.Cells.Item(18, 5).Value2 = 2
$ws.Cells.Item(18, 7).Value2 = 0.182912
$ws.Cells.Item(18, 8).Value2 = 0.365824
$ws.Cells.Item(18, 9).Value2 = 0.129739974361072
$ws.Cells.Item(18, 10).Value2 = 0.09040294133782634
$ws.Cells.Item(18, 11).Value2 = 3
$ws.Cells.Item(18, 13).Value2 = 8.648731
$ws.Cells.Item(18, 14).Value2 = 25.946193
$ws.Cells.Item(18, 15).Value2 = 0.4136938666792861
$ws.Cells.Item(18, 16).Value2 = 0.4684295018482661
$ws.Cells.Item(18, 17).Value2 = 1.581956684672
$ws.Cells.Item(18, 18).Value2 = 9.491740108032001
$ws.Cells.Item(18, 19).Value2 = 0.0536726316563033
$ws.Cells.Item(18, 20).Value2 = 0.04234740477649602

$ws.Cells.Item(19, 5).Value2 = 2
$ws.Cells.Item(19, 7).Value2 = 0.182912
$ws.Cells.Item(19, 8).Value2 = 0.365824
$ws.Cells.Item(19, 9).Value2 = 0.129739974361072
$ws.Cells.Item(19, 10).Value2 = 0.09040294133782634
$ws.Cells.Item(19, 11).Value2 = 3
$ws.Cells.Item(19, 13).Value2 = 1.069144
$ws.Cells.Item(19, 14).Value2 = 3.207432
$ws.Cells.Item(19, 15).Value2 = 0.05114025576665047
$ws.Cells.Item(19, 16).Value2 = 0.05790659824245461
$ws.Cells.Item(19, 17).Value2 = 0.195559267328
$ws.Cells.Item(19, 18).Value2 = 1.173355603968
$ws.Cells.Item(19, 19).Value2 = 0.006634935471983895
$ws.Cells.Item(19, 20).Value2 = 0.005234926803985702

$ws.Cells.Item(20, 5).Value2 = 2
$ws.Cells.Item(20, 7).Value2 = 0.182912
$ws.Cells.Item(20, 8).Value2 = 0.365824
$ws.Cells.Item(20, 9).Value2 = 0.129739974361072
$ws.Cells.Item(20, 10).Value2 = 0.09040294133782634
$ws.Cells.Item(20, 11).Value2 = 3
$ws.Cells.Item(20, 13).Value2 = 1.940482
$ws.Cells.Item(20, 14).Value2 = 5.821446
$ws.Cells.Item(20, 15).Value2 = 0.09281887733605711
$ws.Cells.Item(20, 16).Value2 = 0.1050996980488267
$ws.Cells.Item(20, 17).Value2 = 0.354937443584
$ws.Cells.Item(20, 18).Value2 = 2.129624661504
$ws.Cells.Item(20, 19).Value2 = 0.01204231876580353
$ws.Cells.Item(20, 20).Value2 = 0.009501321837331343

$ws.Cells.Item(21, 5).Value2 = 2
$ws.Cells.Item(21, 7).Value2 = 0.182912
$ws.Cells.Item(21, 8).Value2 = 0.365824
$ws.Cells.Item(21, 9).Value2 = 0.129739974361072
$ws.Cells.Item(21, 10).Value2 = 0.09040294133782634
$ws.Cells.Item(21, 11).Value2 = 2
$ws.Cells.Item(21, 13).Value2 = 7.328590999999999
$ws.Cells.Item(21, 14).Value2 = 14.657182
$ws.Cells.Item(21, 15).Value2 = 0.3505477448773717
$ws.Cells.Item(21, 16).Value2 = 0.2646190314995103
$ws.Cells.Item(21, 17).Value2 = 1.340487236992
$ws.Cells.Item(21, 18).Value2 = 5.361948947968
$ws.Cells.Item(21, 19).Value2 = 0.0454800554327218
$ws.Cells.Item(21, 20).Value2 = 0.02392233878152265
